# Insert a new daily price record as row 30 in the "Pepino ensalada" sheet.
# This pushes the previously-existing rows 30..118 down to 31..119
# (each keeping its own data), and the new row 30 is populated with a
# fresh observation (new date, variety "Alaska", new volume/prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 30:118 down by one row, creating a blank row 30.
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new observation.
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44659
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 100112043
$ws.Range("G30").Value = "Pepino ensalada"
$ws.Range("H30").Value = "Alaska"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 250
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = 10800
$ws.Range("N30").Value = "$/caja 60 unidades"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 180
$ws.Range("Q30").Value = 60
$ws.Range("R30").Value = "Hortaliza"
